# Hortaliza, Mercado Mayorista Lo Valledor de Santiago - Cebollín
# Insert a new week of records (6 rows, fecha 2021-09-22 / serial 44461) right
# before the existing block that starts at row 856 ("Extra" quality, fecha
# 2021-06-10 / serial 44357). This shifts all rows from 856 downward by 6
# (so the old last row 921 becomes row 927), matching the new dimension
# A1:R927.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows right before current row 856 (i.e. at 850:855, since the
# rows that were 850:855 get pushed down to 856:861).
$ws.Rows("850:855").Insert()

$newRows = @(
    @{ Row=850; I="Extra";   O="Provincia de Chacabuco"; D=44461; J=650; K=3000; L=3500; M=3285; P=91 },
    @{ Row=851; I="Extra";   O="Región Metropolitana";   D=44461; J=660; K=3200; L=3500; M=3355; P=93 },
    @{ Row=852; I="Primera"; O="Provincia de Chacabuco"; D=44461; J=790; K=2500; L=2800; M=2667; P=74 },
    @{ Row=853; I="Primera"; O="Región Metropolitana";   D=44461; J=750; K=2500; L=2700; M=2601; P=72 },
    @{ Row=854; I="Segunda"; O="Provincia de Chacabuco"; D=44461; J=250; K=2200; L=2200; M=2200; P=61 },
    @{ Row=855; I="Segunda"; O="Región Metropolitana";   D=44461; J=290; K=2200; L=2200; M=2200; P=61 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2 = 6
    $ws.Cells.Item($row, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($row, 3).Value2 = "Metropolitana"
    $ws.Cells.Item($row, 4).Value2 = $r.D
    $ws.Cells.Item($row, 5).Value2 = 13
    $ws.Cells.Item($row, 6).Value2 = 100112037
    $ws.Cells.Item($row, 7).Value2 = "Cebollín"
    $ws.Cells.Item($row, 8).Value2 = "Sin especificar"
    $ws.Cells.Item($row, 9).Value2 = $r.I
    $ws.Cells.Item($row, 10).Value2 = $r.J
    $ws.Cells.Item($row, 11).Value2 = $r.K
    $ws.Cells.Item($row, 12).Value2 = $r.L
    $ws.Cells.Item($row, 13).Value2 = $r.M
    $ws.Cells.Item($row, 14).Value2 = "`$/paquete 36 unidades"
    $ws.Cells.Item($row, 15).Value2 = $r.O
    $ws.Cells.Item($row, 16).Value2 = $r.P
    $ws.Cells.Item($row, 17).Value2 = 36
    $ws.Cells.Item($row, 18).Value2 = "Hortaliza"
}
